$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in column D (rows 5 and 6)
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 7

# Add new row 9 with data for Aug 8, 2022 (serial date 44781)
$ws.Range("A9").Value = 44781
$ws.Range("B9").Value = 300
$ws.Range("C9").Value = 171
$ws.Range("D9").Value = 43

# Copy the date style from A8 to A9, and the number style from B8:D8 to B9:D9
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("D9").PasteSpecial(-4122)

# Update selection to match target view state
$ws.Range("F6").Select()
